# Update countries & provincias Spain
#
# The underlying "Pais" sheet is a table of COVID-19 stats per country,
# sorted by "Casos totales" (column B) descending. The source data was
# refreshed (new counts for a few countries), which reshuffles sort order
# for some rows. We apply the edit by writing the final, correct values
# directly into the affected cells (letting the engine regenerate the
# shared-strings table itself).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-less direct cell writes: Cells.Item(row, col)
# Columns: A=1 Pais, B=2 Casos totales, C=3 Nuevos casos, D=4 Casos activos,
#          E=5 Recuperados, F=6 Casos criticos, G=7 Muertes hoy, H=8 Muertes

# --- Title / timestamp cell -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 11:34"

# --- Row 18: Belgica gets refreshed totals (no reordering needed) ----------
$ws.Cells.Item(18, 2).Value = 52596
$ws.Cells.Item(18, 3).Value = 585
$ws.Cells.Item(18, 4).Value = 13411
$ws.Cells.Item(18, 5).Value = 30604
$ws.Cells.Item(18, 6).Value = 502
$ws.Cells.Item(18, 7).Value = 60
$ws.Cells.Item(18, 8).Value = 8581

# --- Rows 62-65: Afganistan refreshed and moved above Ghana/Nigeria/Luxemburgo
# Row 62 becomes Afganistan with brand-new figures.
$ws.Cells.Item(62, 1).Value = "Afganistan"
$ws.Cells.Item(62, 2).Value = 4033
$ws.Cells.Item(62, 3).Value = 255
$ws.Cells.Item(62, 4).Value = 502
$ws.Cells.Item(62, 5).Value = 3416
$ws.Cells.Item(62, 6).Value = 7
$ws.Cells.Item(62, 7).Value = 6
$ws.Cells.Item(62, 8).Value = 115

# Row 63 becomes Ghana (previously row 62's data).
$ws.Cells.Item(63, 1).Value = "Ghana"
$ws.Cells.Item(63, 2).Value = 4012
$ws.Cells.Item(63, 3).Value = 0
$ws.Cells.Item(63, 4).Value = 323
$ws.Cells.Item(63, 5).Value = 3671
$ws.Cells.Item(63, 6).Value = 8
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = 18

# Row 64 becomes Nigeria (previously row 63's data).
$ws.Cells.Item(64, 1).Value = "Nigeria"
$ws.Cells.Item(64, 2).Value = 3912
$ws.Cells.Item(64, 3).Value = 0
$ws.Cells.Item(64, 4).Value = 679
$ws.Cells.Item(64, 5).Value = 3116
$ws.Cells.Item(64, 6).Value = 4
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(64, 8).Value = 117

# Row 65 becomes Luxemburgo (previously row 64's data).
$ws.Cells.Item(65, 1).Value = "Luxemburgo"
$ws.Cells.Item(65, 2).Value = 3871
$ws.Cells.Item(65, 3).Value = 0
$ws.Cells.Item(65, 4).Value = 3526
$ws.Cells.Item(65, 5).Value = 245
$ws.Cells.Item(65, 6).Value = 16
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = 100

# --- Rows 173-176: Malaui refreshed and moved above Siria/Macao/Angola -----
# Row 173 becomes Malaui with brand-new figures.
$ws.Cells.Item(173, 1).Value = "Malaui"
$ws.Cells.Item(173, 2).Value = 56
$ws.Cells.Item(173, 3).Value = 13
$ws.Cells.Item(173, 4).Value = 14
$ws.Cells.Item(173, 5).Value = 39
$ws.Cells.Item(173, 6).Value = 1
$ws.Cells.Item(173, 7).Value = 0
$ws.Cells.Item(173, 8).Value = 3

# Row 174 becomes Siria (previously row 173's data).
$ws.Cells.Item(174, 1).Value = "Siria"
$ws.Cells.Item(174, 2).Value = 47
$ws.Cells.Item(174, 3).Value = 0
$ws.Cells.Item(174, 4).Value = 29
$ws.Cells.Item(174, 5).Value = 15
$ws.Cells.Item(174, 6).Value = 0
$ws.Cells.Item(174, 7).Value = 0
$ws.Cells.Item(174, 8).Value = 3

# Row 175 becomes Macao (previously row 174's data).
$ws.Cells.Item(175, 1).Value = "Macao"
$ws.Cells.Item(175, 2).Value = 45
$ws.Cells.Item(175, 3).Value = 0
$ws.Cells.Item(175, 4).Value = 40
$ws.Cells.Item(175, 5).Value = 5
$ws.Cells.Item(175, 6).Value = 1
$ws.Cells.Item(175, 7).Value = 0
$ws.Cells.Item(175, 8).Value = 0

# Row 176 becomes Angola (previously row 175's data).
$ws.Cells.Item(176, 1).Value = "Angola"
$ws.Cells.Item(176, 2).Value = 43
$ws.Cells.Item(176, 3).Value = 0
$ws.Cells.Item(176, 4).Value = 11
$ws.Cells.Item(176, 5).Value = 30
$ws.Cells.Item(176, 6).Value = 0
$ws.Cells.Item(176, 7).Value = 0
$ws.Cells.Item(176, 8).Value = 2

# --- Rows 192-193: Belice and Nueva Caledonia swap order --------------------
$ws.Cells.Item(192, 1).Value = "Belice"
$ws.Cells.Item(192, 2).Value = 18
$ws.Cells.Item(192, 3).Value = 0
$ws.Cells.Item(192, 4).Value = 16
$ws.Cells.Item(192, 5).Value = 0
$ws.Cells.Item(192, 6).Value = 0
$ws.Cells.Item(192, 7).Value = 0
$ws.Cells.Item(192, 8).Value = 2

$ws.Cells.Item(193, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(193, 2).Value = 18
$ws.Cells.Item(193, 3).Value = 0
$ws.Cells.Item(193, 4).Value = 18
$ws.Cells.Item(193, 5).Value = 0
$ws.Cells.Item(193, 6).Value = 0
$ws.Cells.Item(193, 7).Value = 0
$ws.Cells.Item(193, 8).Value = 0

# --- Rows 212-213: Butan and Islas Virgenes Britanicas swap order -----------
$ws.Cells.Item(212, 1).Value = "Butan"
$ws.Cells.Item(212, 2).Value = 7
$ws.Cells.Item(212, 3).Value = 0
$ws.Cells.Item(212, 4).Value = 5
$ws.Cells.Item(212, 5).Value = 2
$ws.Cells.Item(212, 6).Value = 0
$ws.Cells.Item(212, 7).Value = 0
$ws.Cells.Item(212, 8).Value = 0

$ws.Cells.Item(213, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(213, 2).Value = 7
$ws.Cells.Item(213, 3).Value = 0
$ws.Cells.Item(213, 4).Value = 4
$ws.Cells.Item(213, 5).Value = 2
$ws.Cells.Item(213, 6).Value = 0
$ws.Cells.Item(213, 7).Value = 0
$ws.Cells.Item(213, 8).Value = 1
